{"js": "// Append four new paragraphs at the very end of the document body (right\n// before the final section break), matching the diff:\n//\n//   1. \"P\u00f3ki co sygna\u0142y strb i prot nie s\u0105 obs\u0142ugiwane \u2013 znaczy s\u0105 ale nic z nimi nie robi\u0119\"\n//   2. \"Na tym etapie projekt jest syntezowalny\"\n//   3. \"No dobra to teraz podpinanie VIPa\" + \".\"   (two separate runs)\n//   4. \"Wygl\u0105da na to, \u017ce trzeba przepisa\u0107 porty z mojego interfejsu i w topie je po prostu wszystkie podpi\u0105\u0107 pod interfejs z vipa.\"\n//\n// Built and inserted as raw OOXML (Flat-OPC wrapped) via Range.insertOoxml so\n// the third paragraph keeps its text split across two distinct <w:r> runs\n// (as in the source diff) without picking up stray <w:rPr/> padding that a\n// plain insertText()-based run split would leave behind.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p><w:r><w:t>P\u00f3ki co sygna\u0142y strb i prot nie s\u0105 obs\u0142ugiwane \u2013 znaczy s\u0105 ale nic z nimi nie robi\u0119</w:t></w:r></w:p>\" +\n  \"<w:p><w:r><w:t>Na tym etapie projekt jest syntezowalny</w:t></w:r></w:p>\" +\n  \"<w:p><w:r><w:t>No dobra to teraz podpinanie VIPa</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>\" +\n  \"<w:p><w:r><w:t>Wygl\u0105da na to, \u017ce trzeba przepisa\u0107 porty z mojego interfejsu i w topie je po prostu wszystkie podpi\u0105\u0107 pod interfejs z vipa.</w:t></w:r></w:p>\" +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nlastParagraph.getRange(\"After\").insertOoxml(ooxml, \"End\");\n\nawait context.sync();\n", "ps1": "# Append four new paragraphs at the very end of the document body (right\n# before the final section break), matching the diff:\n#\n#   1. \"P\u00f3ki co sygna\u0142y strb i prot nie s\u0105 obs\u0142ugiwane \u2013 znaczy s\u0105 ale nic z nimi nie robi\u0119\"\n#   2. \"Na tym etapie projekt jest syntezowalny\"\n#   3. \"No dobra to teraz podpinanie VIPa\" + \".\"   (two separate runs)\n#   4. \"Wygl\u0105da na to, \u017ce trzeba przepisa\u0107 porty z mojego interfejsu i w topie je po prostu wszystkie podpi\u0105\u0107 pod interfejs z vipa.\"\n#\n# Built and inserted as raw OOXML (Flat-OPC wrapped) via Range.InsertXML so\n# the third paragraph keeps its text split across two distinct <w:r> runs\n# (as in the source diff) without picking up stray <w:rPr/> padding that a\n# plain InsertAfter()/Collapse() split would leave behind.\n\n$d = $word.ActiveDocument\n$r = $d.Content\n$r.Collapse(0)\n\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:r><w:t>P\u00f3ki co sygna\u0142y strb i prot nie s\u0105 obs\u0142ugiwane \u2013 znaczy s\u0105 ale nic z nimi nie robi\u0119</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Na tym etapie projekt jest syntezowalny</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>No dobra to teraz podpinanie VIPa</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Wygl\u0105da na to, \u017ce trzeba przepisa\u0107 porty z mojego interfejsu i w topie je po prostu wszystkie podpi\u0105\u0107 pod interfejs z vipa.</w:t></w:r></w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($xml)\n"}
